$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Re-style the three data tables (slides 14, 15, 16) from the custom
#    "Table_0" style to the built-in "No Style, Table Grid" style.
# ---------------------------------------------------------------------------
$newTableStyleId = "{9E7C24A7-07E3-4946-A8FD-1F118668F11E}"

foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle($newTableStyleId)
        }
    }
}

# ---------------------------------------------------------------------------
# 2) Swap the presentation's theme palette from "Integral / Red Violet" back
#    to the stock "Office" palette (the other theme part - used only by the
#    notes master - keeps the Integral palette; only the colour values move,
#    the font scheme / format scheme were already identical between the two
#    theme parts).
# ---------------------------------------------------------------------------
$cs = $p.SlideMaster.Theme.ThemeColorScheme

$cs.Colors(1).RGB  = 0          # dk1      000000
$cs.Colors(2).RGB  = 16777215   # lt1      FFFFFF
$cs.Colors(3).RGB  = 6968388    # dk2      44546A
$cs.Colors(4).RGB  = 15132391   # lt2      E7E6E6
$cs.Colors(5).RGB  = 13998939   # accent1  5B9BD5
$cs.Colors(6).RGB  = 3243501    # accent2  ED7D31
$cs.Colors(7).RGB  = 10855845   # accent3  A5A5A5
$cs.Colors(8).RGB  = 49407      # accent4  FFC000
$cs.Colors(9).RGB  = 12874308   # accent5  4472C4
$cs.Colors(10).RGB = 4697456    # accent6  70AD47
$cs.Colors(11).RGB = 12673797   # hlink    0563C1
$cs.Colors(12).RGB = 7491477    # folHlink 954F72
